# Insert a new weekly price record for "Ajo" (Chilote variety) at Feria
# Lagunitas de Puerto Montt. The new record is inserted as row 54, which
# pushes all the existing rows 54:151 down by one (to 55:152) — matching
# the diff's shifted dimension (A1:R151 -> A1:R152) and shifted data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 54; everything that was row 54.. becomes 55..
$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with the new record's data.
$ws.Range("A54").Value = 4
$ws.Range("B54").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C54").Value = "Los Lagos"
$ws.Range("D54").Value = 44469
$ws.Range("E54").Value = 10
$ws.Range("F54").Value = 100112003
$ws.Range("G54").Value = "Ajo"
$ws.Range("H54").Value = "Chilote"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 120
$ws.Range("K54").Value = 17000
$ws.Range("L54").Value = 17000
$ws.Range("M54").Value = 17000
$ws.Range("N54").Value = '$/caja 10 kilos'
$ws.Range("O54").Value = "China"
$ws.Range("P54").Value = 1700
$ws.Range("Q54").Value = 10
$ws.Range("R54").Value = "Hortaliza"
